# Update workbook with corrected forecast output:
#  - Rename Sheet1 to "Sales vs PO" and reshape its data (new "Order Week"
#    column inserted, A-column dates shifted forward one week, PO qty zeroed).
#  - Add three new sheets: "Weekly Growth", "Volume Insights", "Prediction Info".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# ---------------------------------------------------------------------------
# Sheet 1: "Sales vs PO"
# ---------------------------------------------------------------------------
# Insert a new column C ("Order Week"); this shifts the old PO_Requested_Qty
# column (with its data + formatting) from C to D automatically.
$ws1.Columns.Item(3).Insert()

# Give the new "Order Week" column (header + 18 data rows) the same
# formatting as the existing "ds" column (header style + date number format),
# then set the header text.
$ws1.Range("A1:A18").Copy()
$ws1.Range("C1:C18").PasteSpecial(-4122)
$ws1.Range("C1").Value = "Order Week"

# Old "ds" values move into the new "Order Week" column, the "ds" column
# itself gets new (one-week-later) dates, and PO_Requested_Qty is reset to 0.
$oldDs = 45537,45544,45551,45558,45565,45572,45579,45586,45593,45600,45607,45614,45621,45628,45635,45642,45649
$newDs = 45543,45550,45557,45564,45571,45578,45585,45592,45599,45606,45613,45620,45627,45634,45641,45648,45655

for ($i = 0; $i -lt 17; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 3).Value = $oldDs[$i]
    $ws1.Cells.Item($row, 1).Value = $newDs[$i]
    $ws1.Cells.Item($row, 4).Value = 0
}

# ---------------------------------------------------------------------------
# Sheet 2: "Weekly Growth"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"

$ws1.Range("A2:A7").Copy()
$ws2.Range("A2:A7").PasteSpecial(-4122)

$wgDs = 45544,45572,45579,45586,45593,45600
$wgPo = 670,340,230,10,680,40
$wgGrowth = 0,-49.25373134328358,-32.35294117647059,-95.65217391304348,6700,-94.11764705882352

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $wgDs[$i]
    $ws2.Cells.Item($row, 2).Value = $wgPo[$i]
    $ws2.Cells.Item($row, 3).Value = $wgGrowth[$i]
}

# ---------------------------------------------------------------------------
# Sheet 3: "Volume Insights"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws1.Range("A1:D1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"

$ws3.Range("A2").Value = 1970
$ws3.Range("B2").Value = 328.3333333333333
$ws3.Range("C2").Value = 680
$ws3.Range("D2").Value = 10

# ---------------------------------------------------------------------------
# Sheet 4: "Prediction Info"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$ws1.Range("A1:A1").Copy()
$ws4.Range("A1:A1").PasteSpecial(-4122)
$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"

$ws4.Range("A2").Value = 93.33333333333326

# Restore the originally active sheet/tab.
$ws1.Activate()
